$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (C) column for all existing data rows (2-116)
#     from 45175 (2023-09-06) to 45177 (2023-09-08)
$ws.Range("C2:C116").Value = 45177

# --- 2. Give row 116 an explicit row height (matches the target workbook)
$ws.Rows.Item(116).RowHeight = 15

# --- 3. Append the two new cutting-notice rows (117 and 118)
$dateFormat = $ws.Cells.Item(116, 2).NumberFormat

function Add-Notice($row, $beteckning, $datum, $forandrad, $area) {
    $ws.Cells.Item($row, 1).Value = $beteckning
    $ws.Cells.Item($row, 2).Value = $datum
    $ws.Cells.Item($row, 2).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 3).Value = $forandrad
    $ws.Cells.Item($row, 3).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 4).Value = "STOCKHOLMS LÄN"
    $ws.Cells.Item($row, 5).Value = "SIGTUNA"
    $ws.Cells.Item($row, 6).Value = "Allmännings- och besparingsskogar"
    $ws.Cells.Item($row, 7).Value = $area
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0
    $ws.Cells.Item($row, 18).WrapText = $true
}

Add-Notice 117 "A 41702-2023" 45176 45177 1.8
Add-Notice 118 "A 41723-2023" 45176 45177 3.7

# Row 117 carries an explicit height like row 116; row 118 keeps the default.
$ws.Rows.Item(117).RowHeight = 15
